$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values in Q2/R2 to whole numbers
$ws.Range("Q2").Value = 805626
$ws.Range("R2").Value = 7353275

# Remove the Starttid (Z2) and Sluttid (AB2) values entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
